# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the leve-profit sheets
# (mirrors the upstream diff: updated currentAveragePrice* / LevePrice* / LeveProfit* columns)

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1121507
$ws.Range("I86").Value = 2315068.8
$ws.Range("J86").Value = 7516.2
$ws.Range("K86").Value = 2315068.8
$ws.Range("L86").Value = 7516.2
$ws.Range("M86").Value = -2313945.8
$ws.Range("N86").Value = -9762.200000000001
$ws.Range("H89").Value = 1121507
$ws.Range("I89").Value = 2315068.8
$ws.Range("J89").Value = 7516.2
$ws.Range("K89").Value = 11575344
$ws.Range("L89").Value = 37581
$ws.Range("M89").Value = -11569728
$ws.Range("N89").Value = -48813
$ws.Range("H100").Value = 2129.6667
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 3000
$ws.Range("N100").Value = -4082
$ws.Range("H111").Value = 3713.923
$ws.Range("I111").Value = 2491.625
$ws.Range("J111").Value = 5669.6
$ws.Range("K111").Value = 7474.875
$ws.Range("L111").Value = 17008.8
$ws.Range("M111").Value = -4407.875
$ws.Range("N111").Value = -23142.8
$ws.Range("H112").Value = 2685.4119
$ws.Range("I112").Value = 489.5
$ws.Range("J112").Value = 2775.0408
$ws.Range("K112").Value = 1468.5
$ws.Range("L112").Value = 8325.1224
$ws.Range("M112").Value = -360.5
$ws.Range("N112").Value = -10541.1224
$ws.Range("H130").Value = 120603
$ws.Range("J130").Value = 120603
$ws.Range("L130").Value = 120603
$ws.Range("N130").Value = -130643
$ws.Range("H132").Value = 13711.493
$ws.Range("I132").Value = 2532.2354
$ws.Range("K132").Value = 7596.706200000001
$ws.Range("M132").Value = -5066.706200000001
$ws.Range("H133").Value = 94352
$ws.Range("J133").Value = 94352
$ws.Range("L133").Value = 94352
$ws.Range("N133").Value = -104472
$ws.Range("H138").Value = 4747.4414
$ws.Range("J138").Value = 5625.981
$ws.Range("L138").Value = 16877.943
$ws.Range("N138").Value = -27157.943
$ws.Range("H141").Value = 6116.636
$ws.Range("J141").Value = 3971
$ws.Range("L141").Value = 11913
$ws.Range("N141").Value = -22273

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1346263.5
$ws.Range("J2").Value = 2928.8
$ws.Range("L2").Value = 2928.8
$ws.Range("N2").Value = -3154.8
$ws.Range("H74").Value = 906.2
$ws.Range("I74").Value = 467
$ws.Range("J74").Value = 1311.6154
$ws.Range("K74").Value = 467
$ws.Range("L74").Value = 1311.6154
$ws.Range("M74").Value = 407
$ws.Range("N74").Value = -3059.6154
$ws.Range("H77").Value = 906.2
$ws.Range("I77").Value = 467
$ws.Range("J77").Value = 1311.6154
$ws.Range("K77").Value = 2335
$ws.Range("L77").Value = 6558.076999999999
$ws.Range("M77").Value = 2033
$ws.Range("N77").Value = -15294.077
$ws.Range("H102").Value = 490985.06
$ws.Range("I102").Value = 653723.5
$ws.Range("K102").Value = 653723.5
$ws.Range("M102").Value = -652101.5
$ws.Range("H116").Value = 1346263.5
$ws.Range("J116").Value = 2928.8
$ws.Range("L116").Value = 2928.8
$ws.Range("N116").Value = -7516.8
$ws.Range("H132").Value = 22582.482
$ws.Range("I132").Value = 31399.63
$ws.Range("K132").Value = 94198.89
$ws.Range("M132").Value = -91668.89

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1346263.5
$ws.Range("J3").Value = 2928.8
$ws.Range("L3").Value = 2928.8
$ws.Range("N3").Value = -3156.8
$ws.Range("H125").Value = 30000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 30000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 30000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -39840
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 3712.1667
$ws.Range("I134").Value = 1440
$ws.Range("K134").Value = 4320
$ws.Range("M134").Value = -1785

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5408.614
$ws.Range("I31").Value = 874.3570999999999
$ws.Range("K31").Value = 874.3570999999999
$ws.Range("M31").Value = -579.3570999999999
$ws.Range("H34").Value = 5408.614
$ws.Range("I34").Value = 874.3570999999999
$ws.Range("K34").Value = 874.3570999999999
$ws.Range("M34").Value = -672.3570999999999
$ws.Range("H131").Value = 48949
$ws.Range("J131").Value = 48949
$ws.Range("L131").Value = 48949
$ws.Range("N131").Value = -59029
$ws.Range("H134").Value = 2036.5294
$ws.Range("I134").Value = 1673.3
$ws.Range("J134").Value = 2555.4285
$ws.Range("K134").Value = 5019.9
$ws.Range("L134").Value = 7666.2855
$ws.Range("M134").Value = -2484.9
$ws.Range("N134").Value = -12736.2855

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 50006388
$ws.Range("I56").Value = 50006388
$ws.Range("K56").Value = 50006388
$ws.Range("M56").Value = -50005858
$ws.Range("H121").Value = 1231
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H131").Value = 14006119
$ws.Range("I131").Value = 55556176
$ws.Range("K131").Value = 166668528
$ws.Range("M131").Value = -166663488

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 400
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H80").Value = 1450860.1
$ws.Range("J80").Value = 45709.43
$ws.Range("L80").Value = 45709.43
$ws.Range("N80").Value = -47705.43
$ws.Range("H83").Value = 1450860.1
$ws.Range("J83").Value = 45709.43
$ws.Range("L83").Value = 228547.15
$ws.Range("N83").Value = -238531.15

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5305.05
$ws.Range("J7").Value = 5480.067
$ws.Range("L7").Value = 5480.067
$ws.Range("N7").Value = -5704.067
$ws.Range("H16").Value = 50002036
$ws.Range("I16").Value = 80001530
$ws.Range("K16").Value = 80001530
$ws.Range("M16").Value = -80001360
$ws.Range("H61").Value = 2108.25
$ws.Range("I61").Value = 1520.3182
$ws.Range("K61").Value = 1520.3182
$ws.Range("M61").Value = -1318.3182
$ws.Range("H93").Value = 8338.25
$ws.Range("I93").Value = 4006.889
$ws.Range("J93").Value = 21332.334
$ws.Range("K93").Value = 4006.889
$ws.Range("L93").Value = 21332.334
$ws.Range("M93").Value = -2758.889
$ws.Range("N93").Value = -23828.334
$ws.Range("H113").Value = 2108.25
$ws.Range("I113").Value = 1520.3182
$ws.Range("K113").Value = 1520.3182
$ws.Range("M113").Value = 649.6818000000001
$ws.Range("H126").Value = 5305.05
$ws.Range("J126").Value = 5480.067
$ws.Range("L126").Value = 16440.201
$ws.Range("N126").Value = -21380.201
